$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-12 from 45224 (2023-10-25)
# to 45233 (2023-11-03), keeping the existing date number format/style.
for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 3).Value = 45233
}
